$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add P1 and Q1, continuing the existing sequence (14, 15) and
# copying the bold/centered/bordered formatting used by the rest of row 1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Data rows 2-25: swap the I<->K and M<->O values, and append the two new
# columns P and Q (value 2 for every row).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2 (new)
}
